$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
for ($i = 0; $i -lt 3; $i++) { $ws.Rows.Item(170).Insert() }

$ws.Cells.Item(2, 1).Value = "Última actualización: 12:12:04"
$ws.Cells.Item(3, 1).Value = "Total filas: 168"
$ws.Cells.Item(16, 1).Value = "05:44:02"
$ws.Cells.Item(16, 2).Value = "06:40"
$ws.Cells.Item(16, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(16, 4).Value = 56
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "06:38:54"
$ws.Cells.Item(17, 2).Value = "06:40"
$ws.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(28, 1).Value = "07:15:48"
$ws.Cells.Item(28, 2).Value = "07:21"
$ws.Cells.Item(28, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(28, 4).Value = 6
$ws.Cells.Item(28, 5).Value = "LP1912"
$ws.Cells.Item(29, 1).Value = "06:56:24"
$ws.Cells.Item(29, 2).Value = "07:21"
$ws.Cells.Item(29, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(29, 4).Value = 25
$ws.Cells.Item(29, 5).Value = "LP1912"
$ws.Cells.Item(51, 1).Value = "08:30:14"
$ws.Cells.Item(51, 2).Value = "08:30"
$ws.Cells.Item(51, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = "LP1912"
$ws.Cells.Item(52, 1).Value = "08:30:14"
$ws.Cells.Item(52, 2).Value = "08:30"
$ws.Cells.Item(52, 3).Value = "15_ABASTO"
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = "LP1912"
$ws.Cells.Item(53, 1).Value = "08:40:59"
$ws.Cells.Item(53, 2).Value = "08:40"
$ws.Cells.Item(53, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = "LP1912"
$ws.Cells.Item(54, 1).Value = "08:40:59"
$ws.Cells.Item(54, 2).Value = "08:40"
$ws.Cells.Item(54, 3).Value = "10_OLMOS"
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = "LP1912"
$ws.Cells.Item(55, 1).Value = "07:52:32"
$ws.Cells.Item(55, 2).Value = "08:41"
$ws.Cells.Item(55, 3).Value = "10_OLMOS"
$ws.Cells.Item(55, 4).Value = 49
$ws.Cells.Item(55, 5).Value = "LP1912"
$ws.Cells.Item(56, 1).Value = "08:30:14"
$ws.Cells.Item(56, 2).Value = "08:41"
$ws.Cells.Item(56, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(56, 4).Value = 11
$ws.Cells.Item(56, 5).Value = "LP1912"
$ws.Cells.Item(64, 1).Value = "08:30:14"
$ws.Cells.Item(64, 2).Value = "08:54"
$ws.Cells.Item(64, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(64, 4).Value = 24
$ws.Cells.Item(64, 5).Value = "LP1912"
$ws.Cells.Item(65, 1).Value = "08:30:14"
$ws.Cells.Item(65, 2).Value = "08:54"
$ws.Cells.Item(65, 3).Value = "215B_EL PATO"
$ws.Cells.Item(65, 4).Value = 24
$ws.Cells.Item(65, 5).Value = "LP1912"
$ws.Cells.Item(75, 1).Value = "08:40:59"
$ws.Cells.Item(75, 2).Value = "09:17"
$ws.Cells.Item(75, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(75, 4).Value = 37
$ws.Cells.Item(75, 5).Value = "LP1912"
$ws.Cells.Item(76, 1).Value = "08:30:14"
$ws.Cells.Item(76, 2).Value = "09:17"
$ws.Cells.Item(76, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(76, 4).Value = 47
$ws.Cells.Item(76, 5).Value = "LP1912"
$ws.Cells.Item(118, 1).Value = "11:13:01"
$ws.Cells.Item(118, 2).Value = "11:15"
$ws.Cells.Item(118, 3).Value = "14_ABASTO"
$ws.Cells.Item(118, 4).Value = 2
$ws.Cells.Item(118, 5).Value = "LP1912"
$ws.Cells.Item(119, 1).Value = "11:13:01"
$ws.Cells.Item(119, 2).Value = "11:15"
$ws.Cells.Item(119, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(119, 4).Value = 2
$ws.Cells.Item(119, 5).Value = "LP1912"
$ws.Cells.Item(128, 1).Value = "10:56:01"
$ws.Cells.Item(128, 2).Value = "11:45"
$ws.Cells.Item(128, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(128, 4).Value = 49
$ws.Cells.Item(128, 5).Value = "LP1912"
$ws.Cells.Item(129, 1).Value = "11:35:40"
$ws.Cells.Item(129, 2).Value = "11:45"
$ws.Cells.Item(129, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(129, 4).Value = 10
$ws.Cells.Item(129, 5).Value = "LP1912"
$ws.Cells.Item(132, 1).Value = "10:07:51"
$ws.Cells.Item(132, 2).Value = "11:52"
$ws.Cells.Item(132, 3).Value = "225_GOMEZ"
$ws.Cells.Item(132, 4).Value = 105
$ws.Cells.Item(132, 5).Value = "LP1912"
$ws.Cells.Item(133, 1).Value = "11:48:20"
$ws.Cells.Item(133, 2).Value = "11:52"
$ws.Cells.Item(133, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(133, 4).Value = 4
$ws.Cells.Item(133, 5).Value = "LP1912"
$ws.Cells.Item(134, 1).Value = "11:35:40"
$ws.Cells.Item(134, 2).Value = "11:53"
$ws.Cells.Item(134, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(134, 4).Value = 18
$ws.Cells.Item(134, 5).Value = "LP1912"
$ws.Cells.Item(135, 1).Value = "11:48:20"
$ws.Cells.Item(135, 2).Value = "11:53"
$ws.Cells.Item(135, 3).Value = "225_GOMEZ"
$ws.Cells.Item(135, 4).Value = 5
$ws.Cells.Item(135, 5).Value = "LP1912"
$ws.Cells.Item(140, 1).Value = "11:55:01"
$ws.Cells.Item(140, 2).Value = "12:10"
$ws.Cells.Item(140, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(140, 4).Value = 15
$ws.Cells.Item(140, 5).Value = "LP1912"
$ws.Cells.Item(141, 1).Value = "11:55:01"
$ws.Cells.Item(141, 2).Value = "12:10"
$ws.Cells.Item(141, 3).Value = "15_ABASTO"
$ws.Cells.Item(141, 4).Value = 15
$ws.Cells.Item(141, 5).Value = "LP1912"
$ws.Cells.Item(142, 1).Value = "12:12:04"
$ws.Cells.Item(142, 2).Value = "12:12"
$ws.Cells.Item(142, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 5).Value = "LP1912"
$ws.Cells.Item(143, 1).Value = "12:12:04"
$ws.Cells.Item(143, 2).Value = "12:13"
$ws.Cells.Item(143, 3).Value = "15_ABASTO"
$ws.Cells.Item(143, 4).Value = 1
$ws.Cells.Item(143, 5).Value = "LP1912"
$ws.Cells.Item(144, 1).Value = "12:12:04"
$ws.Cells.Item(144, 2).Value = "12:17"
$ws.Cells.Item(144, 3).Value = "10_OLMOS"
$ws.Cells.Item(144, 4).Value = 5
$ws.Cells.Item(144, 5).Value = "LP1912"
$ws.Cells.Item(145, 1).Value = "11:35:40"
$ws.Cells.Item(145, 2).Value = "12:21"
$ws.Cells.Item(145, 3).Value = "215C_EL PATO"
$ws.Cells.Item(145, 4).Value = 46
$ws.Cells.Item(145, 5).Value = "LP1912"
$ws.Cells.Item(146, 1).Value = "12:12:04"
$ws.Cells.Item(146, 2).Value = "12:22"
$ws.Cells.Item(146, 3).Value = "215C_EL PATO"
$ws.Cells.Item(146, 4).Value = 10
$ws.Cells.Item(146, 5).Value = "LP1912"
$ws.Cells.Item(147, 1).Value = "11:13:01"
$ws.Cells.Item(147, 2).Value = "12:31"
$ws.Cells.Item(147, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(147, 4).Value = 78
$ws.Cells.Item(147, 5).Value = "LP1912"
$ws.Cells.Item(148, 1).Value = "11:35:40"
$ws.Cells.Item(148, 2).Value = "12:31"
$ws.Cells.Item(148, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(148, 4).Value = 56
$ws.Cells.Item(148, 5).Value = "LP1912"
$ws.Cells.Item(149, 1).Value = "11:48:20"
$ws.Cells.Item(149, 2).Value = "12:32"
$ws.Cells.Item(149, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(149, 4).Value = 44
$ws.Cells.Item(149, 5).Value = "LP1912"
$ws.Cells.Item(150, 1).Value = "11:48:20"
$ws.Cells.Item(150, 2).Value = "12:32"
$ws.Cells.Item(150, 3).Value = "14_ABASTO"
$ws.Cells.Item(150, 4).Value = 44
$ws.Cells.Item(150, 5).Value = "LP1912"
$ws.Cells.Item(151, 1).Value = "12:12:04"
$ws.Cells.Item(151, 2).Value = "12:33"
$ws.Cells.Item(151, 3).Value = "14_ABASTO"
$ws.Cells.Item(151, 4).Value = 21
$ws.Cells.Item(151, 5).Value = "LP1912"
$ws.Cells.Item(152, 1).Value = "11:13:01"
$ws.Cells.Item(152, 2).Value = "12:33"
$ws.Cells.Item(152, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(152, 4).Value = 80
$ws.Cells.Item(152, 5).Value = "LP1912"
$ws.Cells.Item(153, 1).Value = "12:12:04"
$ws.Cells.Item(153, 2).Value = "12:34"
$ws.Cells.Item(153, 3).Value = "15_ABASTO"
$ws.Cells.Item(153, 4).Value = 22
$ws.Cells.Item(153, 5).Value = "LP1912"
$ws.Cells.Item(154, 1).Value = "10:56:01"
$ws.Cells.Item(154, 2).Value = "12:34"
$ws.Cells.Item(154, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(154, 4).Value = 98
$ws.Cells.Item(154, 5).Value = "LP1912"
$ws.Cells.Item(155, 1).Value = "11:35:40"
$ws.Cells.Item(155, 2).Value = "12:36"
$ws.Cells.Item(155, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(155, 4).Value = 61
$ws.Cells.Item(155, 5).Value = "LP1912"
$ws.Cells.Item(156, 1).Value = "12:12:04"
$ws.Cells.Item(156, 2).Value = "12:37"
$ws.Cells.Item(156, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(156, 4).Value = 25
$ws.Cells.Item(156, 5).Value = "LP1912"
$ws.Cells.Item(157, 1).Value = "12:12:04"
$ws.Cells.Item(157, 2).Value = "12:37"
$ws.Cells.Item(157, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(157, 4).Value = 25
$ws.Cells.Item(157, 5).Value = "LP1912"
$ws.Cells.Item(158, 1).Value = "11:48:20"
$ws.Cells.Item(158, 2).Value = "12:47"
$ws.Cells.Item(158, 3).Value = "14_ABASTO"
$ws.Cells.Item(158, 4).Value = 59
$ws.Cells.Item(158, 5).Value = "LP1912"
$ws.Cells.Item(159, 1).Value = "12:12:04"
$ws.Cells.Item(159, 2).Value = "12:48"
$ws.Cells.Item(159, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(159, 4).Value = 36
$ws.Cells.Item(159, 5).Value = "LP1912"
$ws.Cells.Item(160, 1).Value = "12:12:04"
$ws.Cells.Item(160, 2).Value = "12:48"
$ws.Cells.Item(160, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(160, 4).Value = 36
$ws.Cells.Item(160, 5).Value = "LP1912"
$ws.Cells.Item(161, 1).Value = "11:55:01"
$ws.Cells.Item(161, 2).Value = "12:48"
$ws.Cells.Item(161, 3).Value = "14_ABASTO"
$ws.Cells.Item(161, 4).Value = 53
$ws.Cells.Item(161, 5).Value = "LP1912"
$ws.Cells.Item(162, 1).Value = "11:48:20"
$ws.Cells.Item(162, 2).Value = "13:02"
$ws.Cells.Item(162, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(162, 4).Value = 74
$ws.Cells.Item(162, 5).Value = "LP1912"
$ws.Cells.Item(163, 1).Value = "12:12:04"
$ws.Cells.Item(163, 2).Value = "13:03"
$ws.Cells.Item(163, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(163, 4).Value = 51
$ws.Cells.Item(163, 5).Value = "LP1912"
$ws.Cells.Item(164, 1).Value = "11:35:40"
$ws.Cells.Item(164, 2).Value = "13:03"
$ws.Cells.Item(164, 3).Value = "215C_EL PATO"
$ws.Cells.Item(164, 4).Value = 88
$ws.Cells.Item(164, 5).Value = "LP1912"
$ws.Cells.Item(165, 1).Value = "12:12:04"
$ws.Cells.Item(165, 2).Value = "13:04"
$ws.Cells.Item(165, 3).Value = "215C_EL PATO"
$ws.Cells.Item(165, 4).Value = 52
$ws.Cells.Item(165, 5).Value = "LP1912"
$ws.Cells.Item(166, 1).Value = "11:55:01"
$ws.Cells.Item(166, 2).Value = "13:13"
$ws.Cells.Item(166, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(166, 4).Value = 78
$ws.Cells.Item(166, 5).Value = "LP1912"
$ws.Cells.Item(167, 1).Value = "12:12:04"
$ws.Cells.Item(167, 2).Value = "13:17"
$ws.Cells.Item(167, 3).Value = "10_OLMOS"
$ws.Cells.Item(167, 4).Value = 65
$ws.Cells.Item(167, 5).Value = "LP1912"
$ws.Cells.Item(168, 1).Value = "12:12:04"
$ws.Cells.Item(168, 2).Value = "13:24"
$ws.Cells.Item(168, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(168, 4).Value = 72
$ws.Cells.Item(168, 5).Value = "LP1912"
$ws.Cells.Item(169, 1).Value = "12:12:04"
$ws.Cells.Item(169, 2).Value = "13:25"
$ws.Cells.Item(169, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(169, 4).Value = 73
$ws.Cells.Item(169, 5).Value = "LP1912"
$ws.Cells.Item(170, 1).Value = "12:12:04"
$ws.Cells.Item(170, 2).Value = "13:33"
$ws.Cells.Item(170, 3).Value = "14_ABASTO"
$ws.Cells.Item(170, 4).Value = 81
$ws.Cells.Item(170, 5).Value = "LP1912"
$ws.Cells.Item(171, 1).Value = "12:12:04"
$ws.Cells.Item(171, 2).Value = "13:33"
$ws.Cells.Item(171, 3).Value = "215A_EL PATO"
$ws.Cells.Item(171, 4).Value = 81
$ws.Cells.Item(171, 5).Value = "LP1912"
$ws.Cells.Item(172, 1).Value = "12:12:04"
$ws.Cells.Item(172, 2).Value = "13:47"
$ws.Cells.Item(172, 3).Value = "225_GOMEZ"
$ws.Cells.Item(172, 4).Value = 95
$ws.Cells.Item(172, 5).Value = "LP1912"
$ws.Cells.Item(173, 1).Value = "11:55:01"
$ws.Cells.Item(173, 2).Value = "13:49"
$ws.Cells.Item(173, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(173, 4).Value = 114
$ws.Cells.Item(173, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Cells.Item(2, 1).Value = "Última actualización: 12:12:04"
$ws.Cells.Item(26, 1).Value = "12:12:04"
$ws.Cells.Item(26, 2).Value = "12:22"
$ws.Cells.Item(26, 3).Value = "215C_EL PATO"
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = "LP1912"
$ws.Cells.Item(28, 1).Value = "12:12:04"
$ws.Cells.Item(28, 2).Value = "13:04"
$ws.Cells.Item(28, 3).Value = "215C_EL PATO"
$ws.Cells.Item(28, 4).Value = 52
$ws.Cells.Item(28, 5).Value = "LP1912"
$ws.Cells.Item(29, 1).Value = "12:12:04"
$ws.Cells.Item(29, 2).Value = "13:33"
$ws.Cells.Item(29, 3).Value = "215A_EL PATO"
$ws.Cells.Item(29, 4).Value = 81
$ws.Cells.Item(29, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Cells.Item(2, 1).Value = "Última actualización: 12:12:04"
$ws.Cells.Item(3, 1).Value = "Total filas: 23"
$ws.Cells.Item(19, 1).Value = "08:52:33"
$ws.Cells.Item(19, 2).Value = "10:30"
$ws.Cells.Item(19, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(19, 4).Value = 98
$ws.Cells.Item(19, 5).Value = "L6173"
$ws.Cells.Item(20, 1).Value = "10:07:51"
$ws.Cells.Item(20, 2).Value = "10:30"
$ws.Cells.Item(20, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(20, 4).Value = 23
$ws.Cells.Item(20, 5).Value = "L6173"
$ws.Cells.Item(25, 1).Value = "12:12:04"
$ws.Cells.Item(25, 2).Value = "13:12"
$ws.Cells.Item(25, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(25, 4).Value = 60
$ws.Cells.Item(25, 5).Value = "L6203"
$ws.Cells.Item(27, 1).Value = "12:12:04"
$ws.Cells.Item(27, 2).Value = "13:21"
$ws.Cells.Item(27, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(27, 4).Value = 69
$ws.Cells.Item(27, 5).Value = "L6173"
$ws.Cells.Item(28, 1).Value = "12:12:04"
$ws.Cells.Item(28, 2).Value = "13:57"
$ws.Cells.Item(28, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(28, 4).Value = 105
$ws.Cells.Item(28, 5).Value = "L6203"

